# Append four new match rows (rows 3-6) to the "Chris Jordan " sheet,
# matching the existing header/row-2 layout:
#   A dateOfMatch | B venueOfMatch | C matchResult | D ownTeam
#   E opponentTeam | F playerName | G runs | H balls | I numberOf4
#   J numberOf6 | K sr

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(" Oct 18 2020", " Dubai (DSC)", "Match tied (Kings XI won the one-over eliminator)", "Kings XI Punjab", "Mumbai Indians", "Chris Jordan ", "13", "8", "2", "0", "162.50"),
    @(" Oct 24 2020", " Dubai (DSC)", "Kings XI won by 12 runs", "Kings XI Punjab", "Sunrisers Hyderabad", "Chris Jordan ", "7", "12", "0", "0", "58.33"),
    @(" Oct 10 2020", " Abu Dhabi", "KKR won by 2 runs", "Kings XI Punjab", "Kolkata Knight Riders", "Chris Jordan ", "0", "0", "0", "0", "-"),
    @(" Sep 20 2020", " Dubai (DSC)", "Match tied (Capitals won the one-over eliminator)", "Kings XI Punjab", "Delhi Capitals", "Chris Jordan ", "5", "6", "0", "0", "83.33")
)

$startRow = 3
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c - 1]
    }
}
